# Update cryptos list with latest prices and volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column values are stored as literal text (e.g. "27.306.44",
# "1.000", "0.000008643"), so force text format before assigning them
# to stop Excel from auto-converting them into numbers and silently
# reformatting/rounding the displayed text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.306.44"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "1.904.26"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "308.37"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "0.5229"
$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("D8").Value = "0.3784"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "0.07296"
$ws.Range("E9").Value = "  +1.09%  "

$ws.Range("D10").Value = "21.29"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("D11").Value = "0.9045"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").Value = "0.08290"
$ws.Range("E12").Value = "  +8.19%  "

$ws.Range("D13").Value = "97.05"
$ws.Range("E13").Value = "  +2.75%  "

$ws.Range("D14").Value = "1.902.48"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").Value = "5.298"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "0.000008643"
$ws.Range("E17").Value = "  +1.59%  "

$ws.Range("D18").Value = "14.59"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "27.324.16"
$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("D21").Value = "5.099"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("B22").Value = "Cosmos"
$ws.Range("C22").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.449"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("B24").Value = "LidoDAOToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D24").Value = "2.316"
$ws.Range("E24").Value = "  +0.49%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "147.54"
$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "18.26"
$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("D27").Value = "1.750"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "115.50"
$ws.Range("E28").Value = "  +0.91%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "4.854"
$ws.Range("E29").Value = "  +1.19%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "4.927"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.09255"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.05077"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.7992"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.240"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "3.444"
$ws.Range("E35").Value = "  +4.78%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.963"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "2.602"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.5739"
$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02006"
$ws.Range("E39").Value = "  +0.84%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.079"
$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "9.031"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.593"
$ws.Range("E42").Value = "  -0.68%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "116.23"
$ws.Range("E43").Value = "  -2.25%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.1521"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.4881"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.08"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.631"
$ws.Range("E48").Value = "  +2.11%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "38.12"
$ws.Range("E49").Value = "  +1.22%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "63.98"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05944"
$ws.Range("E51").Value = "  +0.31%  "
